$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.903.37'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '2.232.37'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.59'
$ws.Range("E5").Value = '  +1.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.621'
$ws.Range("E6").Value = '  -2.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.72'
$ws.Range("E7").Value = '  -7.00%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.20'
$ws.Range("E10").Value = '  -4.53%  '
$ws.Range("E11").Value = '  +0.70%  '
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").Value = '2.566.64'
$ws.Range("E13").Value = '  -0.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.57'
$ws.Range("E14").Value = '  -3.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.75'
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.66'
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("E17").Value = '  -2.91%  '
$ws.Range("D18").Value = '2.266.50'
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").Value = '41.880.33'
$ws.Range("E19").Value = '  +1.39%  '
$ws.Range("D20").Value = '0.0₃0913'
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.69'
$ws.Range("E21").Value = '  -1.90%  '
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.58'
$ws.Range("E23").Value = '  -2.14%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("E27").Value = '  +0.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.60'
$ws.Range("E28").Value = '  -1.91%  '
$ws.Range("E29").Value = '  -2.77%  '
$ws.Range("E30").Value = '  -2.49%  '
$ws.Range("E31").Value = '  -2.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.53'
$ws.Range("E32").Value = '  -11.08%  '
$ws.Range("E33").Value = '  -1.76%  '
$ws.Range("E34").Value = '  +3.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.71'
$ws.Range("E35").Value = '  +0.54%  '
$ws.Range("E36").Value = '  +4.49%  '
$ws.Range("E37").Value = '  -8.63%  '
$ws.Range("E38").Value = '  -1.96%  '
$ws.Range("E39").Value = '  -5.52%  '
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("E42").Value = '  +1.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.67'
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.53'
$ws.Range("E44").Value = '  -5.85%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.23'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '99.12'
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0968'
$ws.Range("E47").Value = '  +2.52%  '
$ws.Range("D48").Value = '1.473.00'
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.69'
$ws.Range("E49").Value = '  -6.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.28'
$ws.Range("E50").Value = '  +7.79%  '
$ws.Range("E51").Value = '  -2.59%  '
